{"js": "// Change 1: \"...have a chance to return dialogue instead of other actions.\"\n//        -> \"...have a chance to print dialogue.\"\n// Split into separate runs matching the author's incremental edit:\n//   \"() method to have a chance to \" | \"print\" | \" dialogue\" | \".\"\nconst body = context.document.body;\n\nconst oldTailResults = body.search(\"return dialogue instead of other actions.\", { matchCase: true });\noldTailResults.load(\"items\");\nawait context.sync();\n\nconst oldTail = oldTailResults.items[0];\n// Collapse the matched text to empty, leaving an insertion point right\n// after \"() method to have a chance to \".\noldTail.insertText(\"\", Word.InsertLocation.replace);\nawait context.sync();\n\nconst printRun = oldTail.insertText(\"print\", Word.InsertLocation.after);\nawait context.sync();\nconst afterPrint = printRun.getRange(Word.RangeLocation.end);\nawait context.sync();\n\nconst dialogueRun = afterPrint.insertText(\" dialogue\", Word.InsertLocation.after);\nawait context.sync();\nconst afterDialogue = dialogueRun.getRange(Word.RangeLocation.end);\nawait context.sync();\n\nafterDialogue.insertText(\".\", Word.InsertLocation.after);\nawait context.sync();\n\n// Change 2: the \"_GoBack\" bookmark moves from surrounding the \"EatAction\"\n// heading (near the end of the document) to surrounding the sentence\n// \"Before beginning turn, uses map.locationOf(actor).getItems() to find\n// the items available at the location where the zombie is standing on. \"\n// (this is simply Word's automatic tracking of the most-recently-edited\n// range, since the first change above was made after this one).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst sentenceStartResults = body.search(\"Before beginning turn, uses\", { matchCase: true });\nconst sentenceEndResults = body.search(\"the zombie is standing on. \", { matchCase: true });\nsentenceStartResults.load(\"items\");\nsentenceEndResults.load(\"items\");\nawait context.sync();\n\nconst sentenceStart = sentenceStartResults.items[0].getRange(Word.RangeLocation.start);\nconst sentenceEnd = sentenceEndResults.items[0].getRange(Word.RangeLocation.end);\nconst bookmarkTarget = sentenceStart.expandTo(sentenceEnd);\nbookmarkTarget.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Change 1: \"...have a chance to return dialogue instead of other actions.\"\n#        -> \"...have a chance to print dialogue.\"\n# Applied as incremental inserts so the result lands in separate runs,\n# mirroring the author's original typed edit:\n#   \"() method to have a chance to \" | \"print\" | \" dialogue\" | \".\"\n$d = $word.ActiveDocument\n\n$oldTail = $d.Content\n$oldTail.Find.Text = \"return dialogue instead of other actions.\"\n$oldTail.Find.Execute() | Out-Null\n$oldTail.Text = \"\"\n$oldTail.InsertAfter(\"print\")\n$oldTail.Collapse(0)   # wdCollapseEnd\n$oldTail.InsertAfter(\" dialogue\")\n$oldTail.Collapse(0)   # wdCollapseEnd\n$oldTail.InsertAfter(\".\")\n\n# Change 2: the \"_GoBack\" bookmark moves from surrounding the \"EatAction\"\n# heading (near the end of the document) to surrounding the sentence\n# \"Before beginning turn, uses map.locationOf(actor).getItems() to find\n# the items available at the location where the zombie is standing on. \"\n# (Word simply tracks the most-recently-edited range under this name.)\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$startRng = $d.Content\n$startRng.Find.Text = \"Before beginning turn, uses\"\n$startRng.Find.Execute() | Out-Null\n$startPos = $startRng.Start\n\n$endRng = $d.Content\n$endRng.Find.Text = \"the zombie is standing on. \"\n$endRng.Find.Execute() | Out-Null\n$endPos = $endRng.End\n\n$bookmarkRange = $d.Range($startPos, $endPos)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n"}
